$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Tree and Graphs"
$ws.Range("B9").Value = 199
$ws.Range("C9").Value = "Binary Tree Right Side View"

$ws.Hyperlinks.Add($ws.Range("C9"), "https://leetcode.com/problems/binary-tree-right-side-view/", "", "", "Binary Tree Right Side View")

$ws.Range("A9:B9").Style = $ws.Range("A8:B8").Style
$ws.Range("C9").Style = $ws.Range("C8").Style

$ws.Range("D9").Select()
